$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" — mirror the existing header styling (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from an existing styled header cell (H1) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data columns I (I0) and J (IF)
$values = @{
    2 = @(7, 8)
    3 = @(8, 8)
    4 = @(8, 8)
    5 = @(9, 9)
    6 = @(7, 7)
    7 = @(6, 6)
    8 = @(8, 8)
    9 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
